$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("857:865").Insert()
